# Generate Report for Handback
# -----------------------------------------------------------------------
# Row 8 of both the "zh-cn" and "de-de" worksheets corresponds to the
# bdebe3fd-61c2-418b-bfcc-ba55e74d140e.md source file. A handback was
# processed for this file, but the handback's commit isn't the latest
# one available, so the report now records:
#   - Latest Target File  (I8): hyperlink back to the .md file
#   - Latest Handback File (J8): the .xlf file that was handed back
#   - Latest Handback DateTime (K8): when the handback was processed
#   - Error Detail (P8): a warning that the handback isn't current
# The "Error Detail" column is also widened so the long message is
# readable.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a0712c01ed3b439bd4bbdd63488d4d2db435f29/e2e/bdebe3fd-61c2-418b-bfcc-ba55e74d140e.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36b86d95823095331308989ed019fbb4dacb9966/e2e/bdebe3fd-61c2-418b-bfcc-ba55e74d140e.md"
$warning = "The version of handback file is not the latest, current: " + $currentUrl + ", latest: " + $latestUrl + "."
$mdDisplay = "bdebe3fd-61c2-418b-bfcc-ba55e74d140e.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("I8").Value = $mdDisplay
$ws.Hyperlinks.Add($ws.Range("I8"), $currentUrl, "", "", $mdDisplay)

$ws.Range("J8").Value = "bdebe3fd-61c2-418b-bfcc-ba55e74d140e.7f4868d70c6b4e2d8b1aff9af6c52e4a2de270b5.zh-cn.xlf"
$ws.Range("K8").Value = "2016-08-17 20:44:43"
$ws.Range("P8").Value = $warning

$ws.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("I8").Value = $mdDisplay
$ws2.Hyperlinks.Add($ws2.Range("I8"), $currentUrl, "", "", $mdDisplay)

$ws2.Range("J8").Value = "bdebe3fd-61c2-418b-bfcc-ba55e74d140e.7f4868d70c6b4e2d8b1aff9af6c52e4a2de270b5.de-de.xlf"
$ws2.Range("K8").Value = "2016-08-17 20:44:50"
$ws2.Range("P8").Value = $warning

$ws2.Columns.Item(16).ColumnWidth = 39.17
